# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 6
    4  = 3
    5  = 5
    6  = 1
    7  = 3
    8  = 3
    9  = 1
    10 = 1
    11 = 3
    12 = 0
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 3
    21 = 2
    22 = 1
    23 = 3
    24 = 2
    25 = 4
    26 = 4
    27 = 3
    28 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
